# Horarios actualizados - Linea 141 - 262
# Merge the new 08:21:50 scrape pass into each sheet's data table: update the
# "last updated" / "total rows" header cells, then rewrite the full Hora_Scrap..Parada
# block (A6:E<n>) in arrival-time order, same as the source scraper does on every run.
$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws1 = $wb.Worksheets.Item('LP1912')
$ws1.Range('A2').Value = 'Última actualización: 08:21:50'
$ws1.Range('A3').Value = 'Total filas: 108'

$raw1 = @'
04:18:53|04:45|215A_EL PATO|27|LP1912
04:40:33|04:46|215A_EL PATO|6|LP1912
04:18:53|04:53|11_ETCHEVERRY|35|LP1912
04:54:03|04:55|11_ETCHEVERRY|1|LP1912
04:54:03|05:15|14_ABASTO|21|LP1912
04:18:53|05:16|17_ROMERO|58|LP1912
04:40:33|05:16|14_ABASTO|36|LP1912
05:18:56|05:20|14_ABASTO|2|LP1912
04:18:53|05:21|23_HERNANDEZ|63|LP1912
04:40:33|05:22|23_HERNANDEZ|42|LP1912
04:18:53|05:34|215B_EL PATO|76|LP1912
04:54:03|05:35|215B_EL PATO|41|LP1912
04:18:53|05:37|14_ABASTO|79|LP1912
04:18:53|05:46|15_ABASTO|88|LP1912
05:49:40|05:51|14_ABASTO|2|LP1912
04:40:33|06:04|16_SANTA ANA|84|LP1912
05:49:40|06:06|16_SANTA ANA|17|LP1912
04:18:53|06:07|16_SANTA ANA|109|LP1912
04:18:53|06:11|215A_EL PATO|113|LP1912
04:54:03|06:12|215A_EL PATO|78|LP1912
04:18:53|06:13|225_HARAS DEL SUR|115|LP1912
04:40:33|06:14|225_HARAS DEL SUR|94|LP1912
06:15:04|06:15|225_HARAS DEL SUR|0|LP1912
05:18:56|06:20|26_HERNANDEZ|62|LP1912
04:40:33|06:21|26_HERNANDEZ|101|LP1912
05:18:56|06:26|23_HERNANDEZ|68|LP1912
04:40:33|06:27|23_HERNANDEZ|107|LP1912
06:15:04|06:28|23_HERNANDEZ|13|LP1912
04:40:33|06:29|86_EST CHICA-ESC AGRARIA|109|LP1912
04:54:03|06:30|86_EST CHICA-ESC AGRARIA|96|LP1912
04:40:33|06:31|16_SANTA ANA|111|LP1912
05:18:56|06:43|225_C ROCA-H SUR|85|LP1912
04:54:03|06:44|225_C ROCA-H SUR|110|LP1912
05:18:56|06:46|215C_EL PATO|88|LP1912
06:43:40|06:46|225_C ROCA-H SUR|3|LP1912
04:54:03|06:47|215C_EL PATO|113|LP1912
05:18:56|06:58|10_OLMOS|100|LP1912
06:57:30|06:58|14_ABASTO|1|LP1912
05:18:56|06:59|14_ABASTO|101|LP1912
06:15:04|07:00|14_ABASTO|45|LP1912
06:15:04|07:01|16_SANTA ANA|46|LP1912
05:49:40|07:04|23_HERNANDEZ|75|LP1912
05:18:56|07:04|15_ABASTO|106|LP1912
05:49:40|07:05|15_ABASTO|76|LP1912
05:18:56|07:06|225_GOMEZ|108|LP1912
05:49:40|07:07|225_GOMEZ|78|LP1912
05:18:56|07:11|215A_EL PATO|113|LP1912
06:15:04|07:12|215A_EL PATO|57|LP1912
05:18:56|07:15|11_ETCHEVERRY|117|LP1912
06:15:04|07:16|11_ETCHEVERRY|61|LP1912
06:43:40|07:16|16_SANTA ANA|33|LP1912
07:20:40|07:20|10_OLMOS|0|LP1912
05:49:40|07:21|26_HERNANDEZ|92|LP1912
06:15:04|07:23|10_OLMOS|68|LP1912
05:49:40|07:29|10_OLMOS|100|LP1912
05:49:40|07:31|11_ETCHEVERRY|102|LP1912
05:49:40|07:32|84_COLONIA URQUIZA-ESC 49|103|LP1912
06:15:04|07:32|11_ETCHEVERRY|77|LP1912
07:20:40|07:34|23_HERNANDEZ|14|LP1912
05:49:40|07:36|27_EL RETIRO|107|LP1912
06:15:04|07:37|27_EL RETIRO|82|LP1912
05:49:40|07:39|10_OLMOS|110|LP1912
07:20:40|07:46|16_SANTA ANA|26|LP1912
06:43:40|07:47|14_ABASTO|64|LP1912
06:15:04|07:48|14_ABASTO|93|LP1912
06:43:40|07:51|215D_EL PATO|68|LP1912
06:15:04|07:52|215D_EL PATO|97|LP1912
07:47:32|07:55|10_OLMOS|8|LP1912
07:20:40|07:58|16_SANTA ANA|38|LP1912
07:20:40|07:59|23_HERNANDEZ|39|LP1912
06:15:04|08:01|23_HERNANDEZ|106|LP1912
06:43:40|08:03|23_HERNANDEZ|80|LP1912
07:20:40|08:03|11_ETCHEVERRY|43|LP1912
06:57:30|08:06|23_HERNANDEZ|69|LP1912
07:47:32|08:10|16_SANTA ANA|23|LP1912
06:15:04|08:12|15_ABASTO|117|LP1912
07:47:32|08:13|10_OLMOS|26|LP1912
07:47:32|08:16|26_HERNANDEZ|29|LP1912
06:43:40|08:21|26_HERNANDEZ|98|LP1912
06:43:40|08:22|16_P MOR-SANTA ANA|99|LP1912
06:43:40|08:23|215B_EL PATO|100|LP1912
06:43:40|08:27|84_COLONIA URQUIZA-ESC 49|104|LP1912
07:47:32|08:31|23_HERNANDEZ|44|LP1912
07:59:28|08:33|10_OLMOS|34|LP1912
07:59:28|08:34|23_HERNANDEZ|35|LP1912
07:59:28|08:39|26_HERNANDEZ|40|LP1912
06:43:40|08:42|81_EL PELIGRO|119|LP1912
07:20:40|08:43|14_ABASTO|83|LP1912
06:57:30|08:54|17_ROMERO|117|LP1912
07:20:40|09:01|215A_EL PATO|101|LP1912
08:21:50|09:01|23_HERNANDEZ|40|LP1912
07:59:28|09:03|11_ETCHEVERRY|64|LP1912
08:21:50|09:07|26_HERNANDEZ|46|LP1912
07:20:40|09:10|16_P MOR-SANTA ANA|110|LP1912
08:21:50|09:13|10_OLMOS|52|LP1912
07:20:40|09:16|27_EL RETIRO|116|LP1912
08:21:50|09:21|26_HERNANDEZ|60|LP1912
07:59:28|09:22|16_SANTA ANA|83|LP1912
07:47:32|09:22|17_ROMERO|95|LP1912
07:47:32|09:23|11_ETCHEVERRY|96|LP1912
08:21:50|09:29|16_SANTA ANA|68|LP1912
07:47:32|09:32|15_ABASTO|105|LP1912
07:47:32|09:33|10_OLMOS|106|LP1912
08:21:50|09:41|215C_EL PATO|80|LP1912
07:47:32|09:42|215C_EL PATO|115|LP1912
07:47:32|09:43|14_ABASTO|116|LP1912
08:21:50|10:10|16_P MOR-SANTA ANA|109|LP1912
08:21:50|10:12|15_ABASTO|111|LP1912
'@
$lines1 = $raw1 -split "`n"
$count1 = $lines1.Length
$data1 = New-Object 'object[,]' $count1,5
for ($i = 0; $i -lt $count1; $i++) {
    $f = $lines1[$i] -split '\|'
    $data1[$i,0] = $f[0]
    $data1[$i,1] = $f[1]
    $data1[$i,2] = $f[2]
    $data1[$i,3] = [int]$f[3]
    $data1[$i,4] = $f[4]
}
$endRow1 = 5 + $count1
$ws1.Range("A6:E$endRow1").Value = $data1

# ---- Sheet: LP1912-215 ----
$ws2 = $wb.Worksheets.Item('LP1912-215')
$ws2.Range('A2').Value = 'Última actualización: 08:21:50'
$ws2.Range('A3').Value = 'Total filas: 16'

$raw2 = @'
04:18:53|04:45|215A_EL PATO|27|LP1912
04:40:33|04:46|215A_EL PATO|6|LP1912
04:18:53|05:34|215B_EL PATO|76|LP1912
04:54:03|05:35|215B_EL PATO|41|LP1912
04:18:53|06:11|215A_EL PATO|113|LP1912
04:54:03|06:12|215A_EL PATO|78|LP1912
05:18:56|06:46|215C_EL PATO|88|LP1912
04:54:03|06:47|215C_EL PATO|113|LP1912
05:18:56|07:11|215A_EL PATO|113|LP1912
06:15:04|07:12|215A_EL PATO|57|LP1912
06:43:40|07:51|215D_EL PATO|68|LP1912
06:15:04|07:52|215D_EL PATO|97|LP1912
06:43:40|08:23|215B_EL PATO|100|LP1912
07:20:40|09:01|215A_EL PATO|101|LP1912
08:21:50|09:41|215C_EL PATO|80|LP1912
07:47:32|09:42|215C_EL PATO|115|LP1912
'@
$lines2 = $raw2 -split "`n"
$count2 = $lines2.Length
$data2 = New-Object 'object[,]' $count2,5
for ($i = 0; $i -lt $count2; $i++) {
    $f = $lines2[$i] -split '\|'
    $data2[$i,0] = $f[0]
    $data2[$i,1] = $f[1]
    $data2[$i,2] = $f[2]
    $data2[$i,3] = [int]$f[3]
    $data2[$i,4] = $f[4]
}
$endRow2 = 5 + $count2
$ws2.Range("A6:E$endRow2").Value = $data2

# ---- Sheet: 6203-6173 ----
$ws3 = $wb.Worksheets.Item('6203-6173')
$ws3.Range('A2').Value = 'Última actualización: 08:21:50'
$ws3.Range('A3').Value = 'Total filas: 25'

$raw3 = @'
04:18:53|05:43|215A_LA PLATA|85|L6173
04:40:33|05:44|215A_LA PLATA|64|L6173
04:18:53|06:08|215A_LA PLATA|110|L6173
04:40:33|06:09|215A_LA PLATA|89|L6173
05:18:56|06:32|215C_LA PLATA|74|L6203
04:40:33|06:33|215C_LA PLATA|113|L6203
05:18:56|06:59|215B_LP-P MOR-1 Y 57|101|L6173
06:15:04|07:00|215B_LP-P MOR-1 Y 57|45|L6173
06:57:30|07:01|215B_LP-P MOR-1 Y 57|4|L6173
05:49:40|07:07|215B_LP-P MOR-1 Y 57|78|L6173
05:49:40|07:35|215A_LA PLATA|106|L6173
07:20:40|07:37|215A_LA PLATA|17|L6173
06:43:40|08:06|215C_LA PLATA|83|L6203
06:15:04|08:07|215C_LA PLATA|112|L6203
07:20:40|08:09|215C_LA PLATA|49|L6203
06:57:30|08:10|215C_LA PLATA|73|L6203
07:47:32|08:13|215C_LA PLATA|26|L6203
07:59:28|08:17|215C_LA PLATA|18|L6203
06:57:30|08:35|215A_LA PLATA|98|L6173
07:59:28|08:36|215A_LA PLATA|37|L6173
06:43:40|08:38|215A_LA PLATA|115|L6173
08:21:50|08:42|215A_LA PLATA|21|L6173
07:20:40|09:08|215D_LA PLATA|108|L6203
07:47:32|09:09|215D_LA PLATA|82|L6203
08:21:50|10:02|215B_LP-P MOR-40 Y 115|101|L6173
'@
$lines3 = $raw3 -split "`n"
$count3 = $lines3.Length
$data3 = New-Object 'object[,]' $count3,5
for ($i = 0; $i -lt $count3; $i++) {
    $f = $lines3[$i] -split '\|'
    $data3[$i,0] = $f[0]
    $data3[$i,1] = $f[1]
    $data3[$i,2] = $f[2]
    $data3[$i,3] = [int]$f[3]
    $data3[$i,4] = $f[4]
}
$endRow3 = 5 + $count3
$ws3.Range("A6:E$endRow3").Value = $data3

